$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (Price + Volume(1h) columns) for rows 2-51

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.487.13"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.948.97"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "489.97"
$ws.Range("E5").Value = "  -5.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.46"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.945.89"
$ws.Range("E8").Value = "  -5.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.421"
$ws.Range("E9").Value = "  -5.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -6.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("E12").Value = "  -8.69%  "
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.446.75"
$ws.Range("E14").Value = "  -5.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.58"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.470.09"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.947.71"
$ws.Range("E17").Value = "  -5.14%  "
$ws.Range("E18").Value = "  -5.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.64"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.54"
$ws.Range("E21").Value = "  -5.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "316.74"
$ws.Range("E22").Value = "  -7.39%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.460"
$ws.Range("E24").Value = "  -8.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "59.95"
$ws.Range("E25").Value = "  -12.29%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("E29").Value = "  -9.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.45"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.54"
$ws.Range("E31").Value = "  -6.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -5.36%  "
$ws.Range("E33").Value = "  -8.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.43"
$ws.Range("E34").Value = "  -9.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "148.80"
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.36"
$ws.Range("E36").Value = "  -8.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.66"
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0654"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.17"
$ws.Range("E40").Value = "  -7.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.974.08"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "36.18"
$ws.Range("E43").Value = "  -10.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.989"
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.627"
$ws.Range("E45").Value = "  -7.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.52"
$ws.Range("E47").Value = "  -9.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.107.14"
$ws.Range("E48").Value = "  -6.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0234"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.22"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.53"
$ws.Range("E51").Value = "  -10.11%  "
